$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D4","D5","D6","D7","D8","D9","D10","D11","D13","D14","D15","D16","D18","D19","D22","D23","D25","D26","D29","D30","D31","D32","D33","D34","D36","D38","D39","D40","D41","D42","D43","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "30.714.50"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "1.888.40"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "247.88"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.4734"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.2921"
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").Value = "0.06530"
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("D10").Value = "22.03"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").Value = "0.07803"
$ws.Range("E11").Value = "  +1.20%  "
$ws.Range("D12").Value = "1.890.57"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "96.62"
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("D14").Value = "0.7357"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("D15").Value = "5.252"
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("D16").Value = "283.49"
$ws.Range("E16").Value = "  +3.61%  "
$ws.Range("D17").Value = "30.704.47"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "13.24"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").Value = "0.000007533"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "2.139.20"
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").Value = "5.317"
$ws.Range("E22").Value = "  +1.72%  "
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("D25").Value = "9.224"
$ws.Range("D26").Value = "164.64"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").Value = "1.338"
$ws.Range("E29").Value = "  -2.00%  "
$ws.Range("D30").Value = "0.09740"
$ws.Range("E30").Value = "  -2.76%  "
$ws.Range("D31").Value = "1.489"
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("D32").Value = "4.297"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").Value = "4.200"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("D34").Value = "0.04862"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("E35").Value = "  +0.90%  "
$ws.Range("D36").Value = "0.6979"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("D38").Value = "0.01896"
$ws.Range("E38").Value = "  +2.43%  "
$ws.Range("D39").Value = "2.805"
$ws.Range("E39").Value = "  +2.16%  "
$ws.Range("D40").Value = "6.376"
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("D41").Value = "76.08"
$ws.Range("E41").Value = "  +5.44%  "
$ws.Range("D42").Value = "2.003"
$ws.Range("E42").Value = "  +1.89%  "
$ws.Range("D43").Value = "0.4253"
$ws.Range("E43").Value = "  +1.83%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "0.8374"
$ws.Range("E45").Value = "  +0.32%  "
$ws.Range("D46").Value = "101.53"
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("D47").Value = "9.464"
$ws.Range("E47").Value = "  +2.24%  "
$ws.Range("D48").Value = "35.68"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("D49").Value = "7.027"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("D50").Value = "916.26"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").Value = "0.05754"
$ws.Range("E51").Value = "  +2.14%  "

foreach ($addr in $textCells) { $ws.Range($addr).ClearFormats() }
